$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.975.19"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.71%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.423.65"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.34%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.07%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("E5").Value = "'  +1.43%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'144.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +2.61%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.01%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.03%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -0.51%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +1.13%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -0.15%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'4.007.48"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.05%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  -0.56%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'28.31"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +2.03%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.423.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.20%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +0.22%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'61.994.23"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.51%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'6.17"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.16%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'13.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +3.01%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +3.20%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'389.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.79%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'74.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.31%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.65%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +0.06%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +1.00%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.189"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +2.84%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'7.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +3.99%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +0.03%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'8.04"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +1.16%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +0.83%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +3.48%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E33").Value = "'  +1.80%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'5.29"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +6.75%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +0.66%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'168.04"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.29%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'3.456.52"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.15%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +0.87%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'28.36"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +6.26%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -1.20%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.788"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +1.34%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'4.44"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +2.11%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +1.63%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +5.07%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'2.526.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +3.19%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +0.35%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -0.36%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -0.11%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.0265"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +1.56%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -0.92%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +0.03%  "
$ws.Range("E51").Style = "Normal"
